$d = $word.ActiveDocument

# 1. Update the letter date (assign paragraph text directly so the
#    xml:space="preserve" attribute on the run's <w:t> is retained)
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*September 19, 2025*") {
        $p.Range.Text = "September 21, 2025"
        break
    }
}

# 2. Split the mailing address paragraph "4177 Stewart LN, Santa Clara CA 95054"
#    into "4177 Stewart LN" / "Santa Clara, CA 95054" / <blank paragraph>
$addressPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*4177 Stewart LN, Santa Clara CA 95054*") {
        $addressPara = $p
        break
    }
}

if ($addressPara -ne $null) {
    # Replace the single paragraph's text with two paragraphs (one embedded break,
    # no trailing break so the original paragraph mark is reused for the 2nd line).
    $addressPara.Range.Text = "4177 Stewart LN`rSanta Clara, CA 95054"

    # Add a clean blank paragraph right after the new "Santa Clara, CA 95054" line.
    $cityPara = $addressPara.Next()
    $endOfCity = $cityPara.Range
    $endOfCity.Collapse(0)
    $endOfCity.Text = "`r"
}

# 3. Remove the two blank paragraphs that used to follow "Board of Directors"
$boardPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Board of Directors*") {
        $boardPara = $p
        break
    }
}

if ($boardPara -ne $null) {
    $blank1 = $boardPara.Next()
    $blank1.Range.Delete()

    $blank2 = $boardPara.Next()
    $blank2.Range.Delete()
}
